# Inicio da macro SOD
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text cleanup (trailing whitespace removed from labels) ---
$ws.Range("A1").Value = "T11"
$ws.Range("J1").Value = "RG"

# --- New data row 4 (and the mostly-empty row 5 below it) ---
$ws.Range("A4").Value = 20246940865024
$ws.Range("B4").Value = 45511
$ws.Range("C4").Value = 45511
$ws.Range("D4").Value = 20247061008983
$ws.Range("E4").Value = 45519
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = "PRRLOIENTR101"
$ws.Range("H4").Value = 105086940
$ws.Range("I4").Value = "CIT SERVICOS DE RADIOLOGIA LTDA"
$ws.Range("J4").Value = 44630811
$ws.Range("K4").Value = 44702450
$ws.Range("L4").Value = 15508
$ws.Range("M4").Value = 45658
$ws.Range("N4").Value = 45495
$ws.Range("O4").Value = 45495
$ws.Range("P4").Value = 150
$ws.Range("Q4").Value = 0.01
$ws.Range("R4").Value = 0.15
$ws.Range("S4").Value = "CASCAVEL"
$ws.Range("T4").Value = "OES"
$ws.Range("U4").Value = 61
$ws.Range("V4").Value = 395107003
$ws.Range("W4").Value = 395107005
$ws.Range("X4").Value = 395107046
$ws.Range("Y4").Value = "200/5"
$ws.Range("Z4").Value = " "
$ws.Range("AA4").Value = " "
$ws.Range("AB4").Value = " "
$ws.Range("AC4").Value = " "
$ws.Range("A5").Value = " "

# --- Re-apply the column formatting used by the rest of the table ---
# style group 3: general text, word-wrap
# style group 4: dd/mm/yy date, word-wrap
# style group 5: dd/mm/yy date, word-wrap, underlined font
# style group 6: general number, word-wrap, underlined font
$r = $ws.Range("A4")
$r.WrapText = $true

$r = $ws.Range("B4:C4")
$r.WrapText = $true
$r.NumberFormat = "dd/mm/yy"

$r = $ws.Range("D4")
$r.WrapText = $true

$r = $ws.Range("E4")
$r.WrapText = $true
$r.NumberFormat = "dd/mm/yy"

$r = $ws.Range("F4:L4")
$r.WrapText = $true

$r = $ws.Range("M4")
$r.WrapText = $true
$r.NumberFormat = "dd/mm/yy"
$r.Font.Underline = $true

$r = $ws.Range("N4:O4")
$r.WrapText = $true
$r.NumberFormat = "dd/mm/yy"

$r = $ws.Range("P4:R4")
$r.WrapText = $true
$r.NumberFormat = "General"
$r.Font.Underline = $true

$r = $ws.Range("S4:AC4")
$r.WrapText = $true

$r = $ws.Range("A5")
$r.WrapText = $true

# --- Row height for the new data row (auto/wrap height from the source) ---
$ws.Rows.Item(4).RowHeight = 23.85

# --- Update the view / active selection state ---
$ws.Range("H11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
